# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Most cells hold plain text (prices/percentages with formatting like "71.019.66"
# or "  +5.75%  "); a handful of new values (e.g. "6.63", "1.00") would be
# auto-coerced to numbers by Excel, which would silently mangle formatting like
# trailing zeros. Set-TextValue forces those through as literal text while
# restoring the cell style afterward so no visible formatting changes remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '71.110.93'
$ws.Range('E2').Value = '  +5.94%  '
$ws.Range('D3').Value = '3.647.86'
$ws.Range('E3').Value = '  +16.66%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue $ws 'D5' '597.64'
$ws.Range('E5').Value = '  +3.20%  '
Set-TextValue $ws 'D6' '182.73'
$ws.Range('E6').Value = '  +3.06%  '
$ws.Range('D7').Value = '3.644.47'
$ws.Range('E7').Value = '  +16.61%  '
Set-TextValue $ws 'D9' '0.536'
$ws.Range('E9').Value = '  +3.84%  '
$ws.Range('E10').Value = '  +7.04%  '
Set-TextValue $ws 'D11' '6.63'
$ws.Range('E11').Value = '  +3.19%  '
Set-TextValue $ws 'D12' '0.499'
$ws.Range('E12').Value = '  +4.64%  '
Set-TextValue $ws 'D13' '40.92'
$ws.Range('E13').Value = '  +12.23%  '
Set-TextValue $ws 'D14' '0.0000255'
$ws.Range('E14').Value = '  +5.03%  '
$ws.Range('D15').Value = '4.257.29'
$ws.Range('E15').Value = '  +16.72%  '
$ws.Range('D16').Value = '71.041.90'
$ws.Range('E16').Value = '  +5.92%  '
$ws.Range('D17').Value = '3.637.35'
$ws.Range('E17').Value = '  +16.24%  '
$ws.Range('E18').Value = '  +0.92%  '
Set-TextValue $ws 'D19' '7.51'
$ws.Range('E19').Value = '  +6.54%  '
Set-TextValue $ws 'D20' '17.07'
$ws.Range('E20').Value = '  +0.03%  '
Set-TextValue $ws 'D21' '515.07'
$ws.Range('E21').Value = '  +4.99%  '
Set-TextValue $ws 'D22' '9.19'
$ws.Range('E22').Value = '  +17.68%  '
Set-TextValue $ws 'D23' '0.745'
$ws.Range('E23').Value = '  +7.02%  '
Set-TextValue $ws 'D24' '87.63'
$ws.Range('E24').Value = '  +4.39%  '
Set-TextValue $ws 'D25' '2.49'
$ws.Range('E25').Value = '  +9.42%  '
Set-TextValue $ws 'D26' '13.58'
$ws.Range('E26').Value = '  +6.03%  '
Set-TextValue $ws 'D27' '11.03'
$ws.Range('E27').Value = '  +6.83%  '
Set-TextValue $ws 'D28' '1.00'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  +9.78%  '
Set-TextValue $ws 'D30' '8.17'
$ws.Range('E30').Value = '  +1.54%  '
Set-TextValue $ws 'D31' '2.78'
$ws.Range('E31').Value = '  +6.36%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D32' '31.67'
$ws.Range('E32').Value = '  +12.28%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws 'D33' '0.0000111'
$ws.Range('E33').Value = '  +16.92%  '
$ws.Range('E34').Value = '  +3.11%  '
Set-TextValue $ws 'D35' '0.998'
$ws.Range('E35').Value = '  -0.20%  '
Set-TextValue $ws 'D36' '6.14'
$ws.Range('E36').Value = '  +8.24%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws 'D37' '0.346'
$ws.Range('E37').Value = '  +11.32%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws 'D38' '1.01'
$ws.Range('E38').Value = '  +6.25%  '
$ws.Range('E39').Value = '  +7.88%  '
Set-TextValue $ws 'D40' '51.01'
$ws.Range('E40').Value = '  +3.15%  '
$ws.Range('E41').Value = '  +4.67%  '
Set-TextValue $ws 'D42' '45.19'
$ws.Range('E42').Value = '  -6.93%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.132.44'
$ws.Range('E43').Value = '  +11.80%  '
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D44' '8.83'
$ws.Range('E44').Value = '  +5.29%  '
Set-TextValue $ws 'D45' '415.65'
$ws.Range('E45').Value = '  +10.65%  '
Set-TextValue $ws 'D46' '2.81'
$ws.Range('E46').Value = '  +3.81%  '
Set-TextValue $ws 'D47' '28.49'
$ws.Range('E47').Value = '  +13.45%  '
Set-TextValue $ws 'D48' '0.0369'
$ws.Range('E48').Value = '  +5.98%  '
Set-TextValue $ws 'D49' '138.38'
$ws.Range('E49').Value = '  +2.36%  '
Set-TextValue $ws 'D51' '2.49'
$ws.Range('E51').Value = '  +10.87%  '
